$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 127.111115
$ws.Range("I55").Value = 194.25
$ws.Range("J55").Value = 73.40000000000001
$ws.Range("K55").Value = 194.25
$ws.Range("L55").Value = 73.40000000000001
$ws.Range("M55").Value = 19.75
$ws.Range("N55").Value = -501.4

$ws.Range("H62").Value = 6518.6665
$ws.Range("I62").Value = 5470.8184
$ws.Range("J62").Value = 8165.2856
$ws.Range("K62").Value = 5470.8184
$ws.Range("L62").Value = 8165.2856
$ws.Range("M62").Value = -4846.8184
$ws.Range("N62").Value = -9413.285599999999

$ws.Range("H65").Value = 6518.6665
$ws.Range("I65").Value = 5470.8184
$ws.Range("J65").Value = 8165.2856
$ws.Range("K65").Value = 27354.092
$ws.Range("L65").Value = 40826.428
$ws.Range("M65").Value = -24234.092
$ws.Range("N65").Value = -47066.428

$ws.Range("H74").Value = 4900
$ws.Range("I74").Value = 4900
$ws.Range("K74").Value = 4900
$ws.Range("M74").Value = -3964

$ws.Range("H77").Value = 4900
$ws.Range("I77").Value = 4900
$ws.Range("K77").Value = 24500
$ws.Range("M77").Value = -19820

$ws.Range("H86").Value = 66712256
$ws.Range("I86").Value = 13973.5
$ws.Range("J86").Value = 90966180
$ws.Range("K86").Value = 13973.5
$ws.Range("L86").Value = 90966180
$ws.Range("M86").Value = -12850.5
$ws.Range("N86").Value = -90968426

$ws.Range("H89").Value = 66712256
$ws.Range("I89").Value = 13973.5
$ws.Range("J89").Value = 90966180
$ws.Range("K89").Value = 69867.5
$ws.Range("L89").Value = 454830900
$ws.Range("M89").Value = -64251.5
$ws.Range("N89").Value = -454842132

$ws.Range("H96").Value = 553.6667
$ws.Range("I96").Value = 632.8570999999999
$ws.Range("J96").Value = 276.5
$ws.Range("K96").Value = 1898.5713
$ws.Range("L96").Value = 829.5
$ws.Range("M96").Value = -525.5712999999998
$ws.Range("N96").Value = -3575.5

$ws.Range("H133").Value = 93187.55499999999
$ws.Range("J133").Value = 93187.55499999999
$ws.Range("L133").Value = 93187.55499999999
$ws.Range("N133").Value = -103307.555

$ws.Range("H138").Value = 151634.14
$ws.Range("I138").Value = 377665.44
$ws.Range("J138").Value = 5132.3706
$ws.Range("K138").Value = 1132996.32
$ws.Range("L138").Value = 15397.1118
$ws.Range("M138").Value = -1127856.32
$ws.Range("N138").Value = -25677.1118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6467.4116
$ws.Range("I45").Value = 5703.8335
$ws.Range("K45").Value = 5703.8335
$ws.Range("M45").Value = -5326.8335

$ws.Range("H61").Value = 8057.381
$ws.Range("I61").Value = 8680.433999999999
$ws.Range("K61").Value = 8680.433999999999
$ws.Range("M61").Value = -8468.433999999999

$ws.Range("H92").Value = 220029500
$ws.Range("J92").Value = 220029500
$ws.Range("L92").Value = 220029500
$ws.Range("N92").Value = -220034492

$ws.Range("H136").Value = 8057.381
$ws.Range("I136").Value = 8680.433999999999
$ws.Range("K136").Value = 26041.302
$ws.Range("M136").Value = -23491.302

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 56355.668
$ws.Range("I26").Value = 41154.1
$ws.Range("K26").Value = 41154.1
$ws.Range("M26").Value = -40862.1

$ws.Range("H96").Value = 20571.25
$ws.Range("I96").Value = 15623.182
$ws.Range("K96").Value = 15623.182
$ws.Range("M96").Value = -12877.182

$ws.Range("H105").Value = 7585.727
$ws.Range("I105").Value = 7160.3335
$ws.Range("K105").Value = 7160.3335
$ws.Range("M105").Value = -5413.3335

$ws.Range("H134").Value = 3280.8667
$ws.Range("I134").Value = 1468.4445
$ws.Range("K134").Value = 4405.333500000001
$ws.Range("M134").Value = -1870.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3725.4866
$ws.Range("I31").Value = 3186
$ws.Range("K31").Value = 3186
$ws.Range("M31").Value = -2891

$ws.Range("H34").Value = 3725.4866
$ws.Range("I34").Value = 3186
$ws.Range("K34").Value = 3186
$ws.Range("M34").Value = -2984

$ws.Range("H62").Value = 9882.294
$ws.Range("I62").Value = 10111.444
$ws.Range("J62").Value = 9624.5
$ws.Range("K62").Value = 10111.444
$ws.Range("L62").Value = 9624.5
$ws.Range("M62").Value = -9487.444
$ws.Range("N62").Value = -10872.5

$ws.Range("H65").Value = 9882.294
$ws.Range("I65").Value = 10111.444
$ws.Range("J65").Value = 9624.5
$ws.Range("K65").Value = 50557.22
$ws.Range("L65").Value = 48122.5
$ws.Range("M65").Value = -47437.22
$ws.Range("N65").Value = -54362.5

$ws.Range("H99").Value = 222289.95
$ws.Range("I99").Value = 504363
$ws.Range("J99").Value = 5310.6924
$ws.Range("K99").Value = 504363
$ws.Range("L99").Value = 5310.6924
$ws.Range("M99").Value = -502865
$ws.Range("N99").Value = -8306.6924

$ws.Range("H122").Value = 14431.3
$ws.Range("I122").Value = 22049.834
$ws.Range("K122").Value = 66149.50199999999
$ws.Range("M122").Value = -63699.50199999999

$ws.Range("H126").Value = 222289.95
$ws.Range("I126").Value = 504363
$ws.Range("J126").Value = 5310.6924
$ws.Range("K126").Value = 1513089
$ws.Range("L126").Value = 15932.0772
$ws.Range("M126").Value = -1510619
$ws.Range("N126").Value = -20872.0772

$ws.Range("H132").Value = 3996.6775
$ws.Range("I132").Value = 3963.2334
$ws.Range("K132").Value = 11889.7002
$ws.Range("M132").Value = -9359.700199999999

$ws.Range("H134").Value = 1832.8966
$ws.Range("I134").Value = 1201.7084
$ws.Range("J134").Value = 4862.6
$ws.Range("K134").Value = 3605.1252
$ws.Range("L134").Value = 14587.8
$ws.Range("M134").Value = -1070.1252
$ws.Range("N134").Value = -19657.8

$ws.Range("H141").Value = 134777.4
$ws.Range("J141").Value = 134777.4
$ws.Range("L141").Value = 134777.4
$ws.Range("N141").Value = -145137.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 136.75
$ws.Range("I12").Value = 118
$ws.Range("K12").Value = 354
$ws.Range("M12").Value = -181

$ws.Range("H56").Value = 7429.75
$ws.Range("I56").Value = 7429.75
$ws.Range("K56").Value = 7429.75
$ws.Range("M56").Value = -6899.75

$ws.Range("H122").Value = 5273.946
$ws.Range("J122").Value = 6532.5356
$ws.Range("L122").Value = 58792.8204
$ws.Range("N122").Value = -63692.8204

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3785.7144
$ws.Range("I14").Value = 4666.6665
$ws.Range("J14").Value = 3125
$ws.Range("K14").Value = 4666.6665
$ws.Range("L14").Value = 3125
$ws.Range("M14").Value = -4498.6665
$ws.Range("N14").Value = -3461

$ws.Range("H70").Value = 9750.5
$ws.Range("J70").Value = 11110.728
$ws.Range("L70").Value = 11110.728
$ws.Range("N70").Value = -11650.728

$ws.Range("H73").Value = 9750.5
$ws.Range("J73").Value = 11110.728
$ws.Range("L73").Value = 11110.728
$ws.Range("N73").Value = -12982.728

$ws.Range("H126").Value = 17837.428
$ws.Range("I126").Value = 27192
$ws.Range("J126").Value = 12640.444
$ws.Range("K126").Value = 81576
$ws.Range("L126").Value = 37921.33199999999
$ws.Range("M126").Value = -79106
$ws.Range("N126").Value = -42861.33199999999

$ws.Range("H132").Value = 2983.36
$ws.Range("I132").Value = 3091.15
$ws.Range("J132").Value = 2552.2
$ws.Range("K132").Value = 9273.450000000001
$ws.Range("L132").Value = 7656.599999999999
$ws.Range("M132").Value = -6743.450000000001
$ws.Range("N132").Value = -12716.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18627.56
$ws.Range("I7").Value = 22326.033
$ws.Range("J7").Value = 8540.817999999999
$ws.Range("K7").Value = 22326.033
$ws.Range("L7").Value = 8540.817999999999
$ws.Range("M7").Value = -22214.033
$ws.Range("N7").Value = -8764.817999999999

$ws.Range("H16").Value = 2324.8215
$ws.Range("I16").Value = 1824.7084
$ws.Range("J16").Value = 5325.5
$ws.Range("K16").Value = 1824.7084
$ws.Range("L16").Value = 5325.5
$ws.Range("M16").Value = -1654.7084
$ws.Range("N16").Value = -5665.5

$ws.Range("H126").Value = 18627.56
$ws.Range("I126").Value = 22326.033
$ws.Range("J126").Value = 8540.817999999999
$ws.Range("K126").Value = 66978.099
$ws.Range("L126").Value = 25622.454
$ws.Range("M126").Value = -64508.099
$ws.Range("N126").Value = -30562.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7565.6
$ws.Range("J41").Value = 7528.4443
$ws.Range("L41").Value = 7528.4443
$ws.Range("N41").Value = -8308.444299999999

$ws.Range("H45").Value = 17000
$ws.Range("I45").Value = 15000
$ws.Range("K45").Value = 15000
$ws.Range("M45").Value = -14509

$ws.Range("H58").Value = 4264261.5
$ws.Range("I58").Value = 23999
$ws.Range("K58").Value = 23999
$ws.Range("M58").Value = -23691

$ws.Range("H62").Value = 196709.61
$ws.Range("I62").Value = 344797.1
$ws.Range("J62").Value = 11600.25
$ws.Range("K62").Value = 344797.1
$ws.Range("L62").Value = 11600.25
$ws.Range("M62").Value = -344173.1
$ws.Range("N62").Value = -12848.25

$ws.Range("H65").Value = 196709.61
$ws.Range("I65").Value = 344797.1
$ws.Range("J65").Value = 11600.25
$ws.Range("K65").Value = 1723985.5
$ws.Range("L65").Value = 58001.25
$ws.Range("M65").Value = -1720865.5
$ws.Range("N65").Value = -64241.25

$ws.Range("H100").Value = 36164.273
$ws.Range("I100").Value = 9134.333000000001
$ws.Range("J100").Value = 68600.2
$ws.Range("K100").Value = 18268.666
$ws.Range("L100").Value = 137200.4
$ws.Range("M100").Value = -17727.666
$ws.Range("N100").Value = -138282.4

$ws.Range("H122").Value = 6946.9
$ws.Range("I122").Value = 5062.5713
$ws.Range("K122").Value = 15187.7139
$ws.Range("M122").Value = -12737.7139
